# Insert a new data row at sheet row 132 (pushing the existing rows
# 132:165 down to 133:166), then populate the new row:
#   - columns A,B,C,E,F,G,H,I,N,Q,R copied from the row that used to be 132
#     (now sitting at row 133, since it was an exact duplicate of the
#     surrounding metadata)
#   - columns D,J,K,L,M,O,P set to the new record's values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 132..165 down to 133..166, leaving row 132 blank.
$ws.Rows.Item(132).Insert()

# Columns that stay the same as the (now shifted) neighbouring record.
$unchangedCols = @(1, 2, 3, 5, 6, 7, 8, 9, 14, 17, 18)
foreach ($c in $unchangedCols) {
    $ws.Cells.Item(132, $c).Value = $ws.Cells.Item(133, $c).Value2()
}

# New values for the inserted record.
$ws.Cells.Item(132, 4).Value  = 44511            # D - Fecha
$ws.Cells.Item(132, 10).Value = 280              # J - Volumen
$ws.Cells.Item(132, 11).Value = 12000            # K - Precio minimo
$ws.Cells.Item(132, 12).Value = 14000            # L - Precio maximo
$ws.Cells.Item(132, 13).Value = 13071            # M - Precio promedio ponderado
$ws.Cells.Item(132, 15).Value = "Región del Maule"  # O - Origen
$ws.Cells.Item(132, 16).Value = 523              # P - Precio $/Kg
